$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 79.2
$ws.Range("I5").Value = 68.666664
$ws.Range("J5").Value = 95
$ws.Range("K5").Value = 68.666664
$ws.Range("L5").Value = 95
$ws.Range("M5").Value = 46.333336
$ws.Range("N5").Value = -325

$ws.Range("H18").Value = 1192.2222
$ws.Range("I18").Value = 1114.1666
$ws.Range("J18").Value = 1816.6666
$ws.Range("K18").Value = 1114.1666
$ws.Range("L18").Value = 1816.6666
$ws.Range("M18").Value = -830.1666
$ws.Range("N18").Value = -2384.6666

$ws.Range("H40").Value = 1468.9656
$ws.Range("I40").Value = 1378.5714
$ws.Range("J40").Value = 1553.3334
$ws.Range("K40").Value = 1378.5714
$ws.Range("L40").Value = 1553.3334
$ws.Range("M40").Value = -1203.5714
$ws.Range("N40").Value = -1903.3334

$ws.Range("H55").Value = 127.76923
$ws.Range("I55").Value = 32.75
$ws.Range("J55").Value = 170
$ws.Range("K55").Value = 32.75
$ws.Range("L55").Value = 170
$ws.Range("M55").Value = 181.25
$ws.Range("N55").Value = -598

$ws.Range("H58").Value = 2902.4546
$ws.Range("I58").Value = 321.16666
$ws.Range("J58").Value = 6000
$ws.Range("K58").Value = 963.4999799999999
$ws.Range("L58").Value = 18000
$ws.Range("M58").Value = -813.4999799999999
$ws.Range("N58").Value = -18300

$ws.Range("H103").Value = 173.26086
$ws.Range("I103").Value = 167.63158
$ws.Range("K103").Value = 502.8947400000001
$ws.Range("M103").Value = 83.10525999999993

$ws.Range("H131").Value = 1408.8
$ws.Range("I131").Value = 1160.1666
$ws.Range("J131").Value = 2403.3333
$ws.Range("K131").Value = 3480.4998
$ws.Range("L131").Value = 7209.999899999999
$ws.Range("M131").Value = 1559.5002
$ws.Range("N131").Value = -17289.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3747.2036
$ws.Range("I132").Value = 3484.634
$ws.Range("J132").Value = 4575.3076
$ws.Range("K132").Value = 10453.902
$ws.Range("L132").Value = 13725.9228
$ws.Range("M132").Value = -7923.902
$ws.Range("N132").Value = -18785.9228

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2019.3448
$ws.Range("I105").Value = 1908.6666
$ws.Range("K105").Value = 1908.6666
$ws.Range("M105").Value = -161.6666

$ws.Range("H107").Value = 1616.7778
$ws.Range("I107").Value = 956
$ws.Range("K107").Value = 956
$ws.Range("M107").Value = 964

$ws.Range("H134").Value = 19578.018
$ws.Range("I134").Value = 23333.29
$ws.Range("J134").Value = 2679.3
$ws.Range("K134").Value = 69999.87
$ws.Range("L134").Value = 8037.900000000001
$ws.Range("M134").Value = -67464.87
$ws.Range("N134").Value = -13107.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 384.58823
$ws.Range("I22").Value = 363.42856
$ws.Range("J22").Value = 483.33334
$ws.Range("K22").Value = 363.42856
$ws.Range("L22").Value = 483.33334
$ws.Range("M22").Value = -13.42856
$ws.Range("N22").Value = -1183.33334

$ws.Range("H31").Value = 4547581.5
$ws.Range("I31").Value = 2296.16
$ws.Range("J31").Value = 10528221
$ws.Range("K31").Value = 2296.16
$ws.Range("L31").Value = 10528221
$ws.Range("M31").Value = -2001.16
$ws.Range("N31").Value = -10528811

$ws.Range("H34").Value = 4547581.5
$ws.Range("I34").Value = 2296.16
$ws.Range("J34").Value = 10528221
$ws.Range("K34").Value = 2296.16
$ws.Range("L34").Value = 10528221
$ws.Range("M34").Value = -2094.16
$ws.Range("N34").Value = -10528625

$ws.Range("H60").Value = 12718.3
$ws.Range("I60").Value = 4000
$ws.Range("K60").Value = 4000
$ws.Range("M60").Value = -3489

$ws.Range("H62").Value = 37041036
$ws.Range("I62").Value = 4397
$ws.Range("J62").Value = 83336830
$ws.Range("K62").Value = 4397
$ws.Range("L62").Value = 83336830
$ws.Range("M62").Value = -3773
$ws.Range("N62").Value = -83338078

$ws.Range("H65").Value = 37041036
$ws.Range("I65").Value = 4397
$ws.Range("J65").Value = 83336830
$ws.Range("K65").Value = 21985
$ws.Range("L65").Value = 416684150
$ws.Range("M65").Value = -18865
$ws.Range("N65").Value = -416690390

$ws.Range("H122").Value = 1008.6667
$ws.Range("I122").Value = 327
$ws.Range("J122").Value = 1787.7142
$ws.Range("K122").Value = 981
$ws.Range("L122").Value = 5363.142599999999
$ws.Range("M122").Value = 1469
$ws.Range("N122").Value = -10263.1426

$ws.Range("H132").Value = 3059.037
$ws.Range("I132").Value = 2527.6667
$ws.Range("J132").Value = 4121.778
$ws.Range("K132").Value = 7583.000100000001
$ws.Range("L132").Value = 12365.334
$ws.Range("M132").Value = -5053.000100000001
$ws.Range("N132").Value = -17425.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1885.5555
$ws.Range("I34").Value = 1000
$ws.Range("J34").Value = 1996.25
$ws.Range("K34").Value = 3000
$ws.Range("L34").Value = 5988.75
$ws.Range("M34").Value = -2916
$ws.Range("N34").Value = -6156.75

$ws.Range("H40").Value = 51.0625
$ws.Range("I40").Value = 59.615383
$ws.Range("J40").Value = 14
$ws.Range("K40").Value = 238.461532
$ws.Range("L40").Value = 56
$ws.Range("M40").Value = -169.461532
$ws.Range("N40").Value = -194

$ws.Range("H68").Value = 1328.2681
$ws.Range("J68").Value = 1440
$ws.Range("L68").Value = 4320
$ws.Range("N68").Value = -5942

$ws.Range("H71").Value = 1328.2681
$ws.Range("J71").Value = 1440
$ws.Range("L71").Value = 12960
$ws.Range("N71").Value = -21072

$ws.Range("H86").Value = 3087.5
$ws.Range("I86").Value = 350
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 1050
$ws.Range("L86").Value = 12000
$ws.Range("M86").Value = 136
$ws.Range("N86").Value = -14372

$ws.Range("H89").Value = 3087.5
$ws.Range("I89").Value = 350
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 3150
$ws.Range("L89").Value = 36000
$ws.Range("M89").Value = 2778
$ws.Range("N89").Value = -47856

$ws.Range("H97").Value = 308.9
$ws.Range("I97").Value = 286.125
$ws.Range("J97").Value = 400
$ws.Range("K97").Value = 858.375
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = -362.375
$ws.Range("N97").Value = -2192

$ws.Range("H107").Value = 989.2222
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 988.58826
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 2965.76478
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -6805.76478

$ws.Range("H122").Value = 883.4
$ws.Range("J122").Value = 868
$ws.Range("L122").Value = 7812
$ws.Range("N122").Value = -12712

$ws.Range("H125").Value = 3735.4546
$ws.Range("I125").Value = 1863.3334
$ws.Range("J125").Value = 4437.5
$ws.Range("K125").Value = 5590.0002
$ws.Range("L125").Value = 13312.5
$ws.Range("M125").Value = -670.0002000000004
$ws.Range("N125").Value = -23152.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 902.375
$ws.Range("I97").Value = 901.73334
$ws.Range("J97").Value = 903.44446
$ws.Range("K97").Value = 901.73334
$ws.Range("L97").Value = 903.44446
$ws.Range("M97").Value = -405.73334
$ws.Range("N97").Value = -1895.44446

$ws.Range("H132").Value = 76615.63
$ws.Range("I132").Value = 127094.25
$ws.Range("J132").Value = 3192.182
$ws.Range("K132").Value = 381282.75
$ws.Range("L132").Value = 9576.545999999998
$ws.Range("M132").Value = -378752.75
$ws.Range("N132").Value = -14636.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2082.2
$ws.Range("I46").Value = 1736.4546
$ws.Range("J46").Value = 2353.8572
$ws.Range("K46").Value = 1736.4546
$ws.Range("L46").Value = 2353.8572
$ws.Range("M46").Value = -1548.4546
$ws.Range("N46").Value = -2729.8572

$ws.Range("H68").Value = 1040
$ws.Range("I68").Value = 1040
$ws.Range("K68").Value = 1040
$ws.Range("M68").Value = -291

$ws.Range("H71").Value = 1040
$ws.Range("I71").Value = 1040
$ws.Range("K71").Value = 5200
$ws.Range("M71").Value = -1456

$ws.Range("H93").Value = 1288306
$ws.Range("I93").Value = 1803041.6
$ws.Range("J93").Value = 1466.6666
$ws.Range("K93").Value = 1803041.6
$ws.Range("L93").Value = 1466.6666
$ws.Range("M93").Value = -1801793.6
$ws.Range("N93").Value = -3962.6666

$ws.Range("H132").Value = 5836.814
$ws.Range("I132").Value = 9217.044
$ws.Range("J132").Value = 1949.55
$ws.Range("K132").Value = 27651.132
$ws.Range("L132").Value = 5848.65
$ws.Range("M132").Value = -25121.132
$ws.Range("N132").Value = -10908.65

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1336.1842
$ws.Range("I132").Value = 1037.6
$ws.Range("J132").Value = 2455.875
$ws.Range("K132").Value = 3112.8
$ws.Range("L132").Value = 7367.625
$ws.Range("M132").Value = -582.7999999999997
$ws.Range("N132").Value = -12427.625
